# Refreshes crypto price/volume figures and re-ranks GateToken into row 7
# (rows 7-18 shift down to 8-19), per the "Updated symbol list" GitHub Action run.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '''286.95'
$ws.Range("E2").Value = '''2.76%'
$ws.Range("D3").Value = '''28.68'
$ws.Range("E3").Value = '''4.35%'
$ws.Range("D4").Value = '''5.050'
$ws.Range("E4").Value = '''4.49%'
$ws.Range("D5").Value = '''0.06666'
$ws.Range("E5").Value = '''4.52%'
$ws.Range("E6").Value = '''4.43%'
$ws.Range("B7").Value = 'GateToken'
$ws.Range("C7").Value = 'https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt'
$ws.Range("D7").Value = '''3.389'
$ws.Range("E7").Value = '''1.82%'
$ws.Range("B8").Value = 'FTXToken'
$ws.Range("C8").Value = 'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt'
$ws.Range("D8").Value = '''1.374'
$ws.Range("E8").Value = '''4.32%'
$ws.Range("B9").Value = 'MXToken'
$ws.Range("C9").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range("D9").Value = '''0.9403'
$ws.Range("E9").Value = '''4.92%'
$ws.Range("B10").Value = 'WazirX'
$ws.Range("C10").Value = 'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx'
$ws.Range("D10").Value = '''0.1557'
$ws.Range("E10").Value = '''0.77%'
$ws.Range("B11").Value = 'LiechtensteinCryptoassetsExchange'
$ws.Range("C11").Value = 'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx'
$ws.Range("D11").Value = '''0.06546'
$ws.Range("E11").Value = '''1.11%'
$ws.Range("B12").Value = 'MandalaExchangeToken'
$ws.Range("C12").Value = 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'
$ws.Range("D12").Value = '''0.07598'
$ws.Range("E12").Value = '''1.07%'
$ws.Range("B13").Value = 'BitrueCoin'
$ws.Range("C13").Value = 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'
$ws.Range("D13").Value = '''0.02959'
$ws.Range("E13").Value = '''0.63%'
$ws.Range("B14").Value = 'BitMartToken'
$ws.Range("C14").Value = 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'
$ws.Range("D14").Value = '''0.08989'
$ws.Range("E14").Value = '''-0.10%'
$ws.Range("B15").Value = 'BitForexToken'
$ws.Range("C15").Value = 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'
$ws.Range("D15").Value = '''0.001600'
$ws.Range("E15").Value = '''1.71%'
$ws.Range("B16").Value = 'CoinExToken'
$ws.Range("C16").Value = 'https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet'
$ws.Range("D16").Value = '''0.04495'
$ws.Range("E16").Value = '''2.28%'
$ws.Range("B17").Value = 'One'
$ws.Range("C17").Value = 'https://coinranking.com/coin/6Lga5NiXX3rT+one-one'
$ws.Range("D17").Value = '''0.0006465'
$ws.Range("E17").Value = '''0.00%'
$ws.Range("B18").Value = 'TigerCash'
$ws.Range("C18").Value = 'https://coinranking.com/coin/6hIn06L2+tigercash-tch'
$ws.Range("D18").Value = '''0.006352'
$ws.Range("E18").Value = '''4.33%'
$ws.Range("B19").Value = 'LEO'
$ws.Range("C19").Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Range("D19").Value = '''3.445'
$ws.Range("E19").Value = '''-1.17%'
$ws.Range("D20").Value = '''2.252'
$ws.Range("E20").Value = '''0.90%'
$ws.Range("E21").Value = '''2.17%'
$ws.Range("E22").Value = '''-2.91%'
$ws.Range("D23").Value = '''4.087'
$ws.Range("E23").Value = '''4.69%'
$ws.Range("D24").Value = '''0.1554'
$ws.Range("E24").Value = '''3.39%'
$ws.Range("E25").Value = '''0.59%'
$ws.Range("D26").Value = '''0.004497'
$ws.Range("E26").Value = '''5.08%'
$ws.Range("D27").Value = '''0.0001251'
$ws.Range("E27").Value = '''6.05%'
$ws.Range("D28").Value = '''0.0001620'
$ws.Range("E28").Value = '''-2.04%'
$ws.Range("D40").Value = '''0.04202'
$ws.Range("E40").Value = '''3.38%'
$ws.Range("D41").Value = '''0.006749'
$ws.Range("E41").Value = '''1.16%'
$ws.Range("D42").Value = '''0.1255'
$ws.Range("E42").Value = '''-10.74%'
$ws.Range("D43").Value = '''0.002022'
$ws.Range("E43").Value = '''-2.31%'
$ws.Range("D44").Value = '''0.01229'
$ws.Range("E44").Value = '''11.43%'
$ws.Range("D45").Value = '''0.00005669'
$ws.Range("E45").Value = '''2.27%'
$ws.Range("E46").Value = '''25.93%'
$ws.Range("D47").Value = '''0.01308'
$ws.Range("E47").Value = '''-29.25%'
